$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1306.0952
$ws.Range("J28").Value = 950.4
$ws.Range("L28").Value = 950.4
$ws.Range("N28").Value = -1920.4
$ws.Range("H92").Value = 3545.9143
$ws.Range("I92").Value = 3490.48
$ws.Range("J92").Value = 3684.5
$ws.Range("K92").Value = 3490.48
$ws.Range("L92").Value = 3684.5
$ws.Range("M92").Value = -2242.48
$ws.Range("N92").Value = -6180.5
$ws.Range("H98").Value = 29415.117
$ws.Range("J98").Value = 3302.2
$ws.Range("L98").Value = 3302.2
$ws.Range("N98").Value = -6298.2
$ws.Range("H99").Value = 241.09091
$ws.Range("I99").Value = 222.66667
$ws.Range("K99").Value = 668.00001
$ws.Range("M99").Value = 829.99999
$ws.Range("H106").Value = 3202.3684
$ws.Range("I106").Value = 2479
$ws.Range("J106").Value = 5533.222
$ws.Range("K106").Value = 2479
$ws.Range("L106").Value = 5533.222
$ws.Range("M106").Value = -1848
$ws.Range("N106").Value = -6795.222
$ws.Range("H122").Value = 29415.117
$ws.Range("J122").Value = 3302.2
$ws.Range("L122").Value = 9906.599999999999
$ws.Range("N122").Value = -14806.6
$ws.Range("H132").Value = 2465.5083
$ws.Range("I132").Value = 2586.9272
$ws.Range("K132").Value = 7760.7816
$ws.Range("M132").Value = -5230.7816
$ws.Range("H138").Value = 2163.3635
$ws.Range("I138").Value = 1360.9615
$ws.Range("J138").Value = 2882.7585
$ws.Range("K138").Value = 4082.8845
$ws.Range("L138").Value = 8648.2755
$ws.Range("M138").Value = 1057.1155
$ws.Range("N138").Value = -18928.2755

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1023.65955
$ws.Range("I2").Value = 989.09375
$ws.Range("J2").Value = 1097.4
$ws.Range("K2").Value = 989.09375
$ws.Range("L2").Value = 1097.4
$ws.Range("M2").Value = -876.09375
$ws.Range("N2").Value = -1323.4
$ws.Range("H74").Value = 28140.234
$ws.Range("I74").Value = 1455.7368
$ws.Range("K74").Value = 1455.7368
$ws.Range("M74").Value = -581.7367999999999
$ws.Range("H77").Value = 28140.234
$ws.Range("I77").Value = 1455.7368
$ws.Range("K77").Value = 7278.683999999999
$ws.Range("M77").Value = -2910.683999999999
$ws.Range("H116").Value = 1023.65955
$ws.Range("I116").Value = 989.09375
$ws.Range("J116").Value = 1097.4
$ws.Range("K116").Value = 989.09375
$ws.Range("L116").Value = 1097.4
$ws.Range("M116").Value = 1304.90625
$ws.Range("N116").Value = -5685.4
$ws.Range("H132").Value = 2656.02
$ws.Range("I132").Value = 2570
$ws.Range("J132").Value = 4003.6667
$ws.Range("K132").Value = 7710
$ws.Range("L132").Value = 12011.0001
$ws.Range("M132").Value = -5180
$ws.Range("N132").Value = -17071.0001
$ws.Range("H139").Value = 125000
$ws.Range("J139").Value = 125000
$ws.Range("L139").Value = 125000
$ws.Range("N139").Value = -135280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1023.65955
$ws.Range("I3").Value = 989.09375
$ws.Range("J3").Value = 1097.4
$ws.Range("K3").Value = 989.09375
$ws.Range("L3").Value = 1097.4
$ws.Range("M3").Value = -875.09375
$ws.Range("N3").Value = -1325.4
$ws.Range("H86").Value = 21278756
$ws.Range("I86").Value = 33335510
$ws.Range("J86").Value = 2127.7058
$ws.Range("K86").Value = 33335510
$ws.Range("L86").Value = 2127.7058
$ws.Range("M86").Value = -33334387
$ws.Range("N86").Value = -4373.7058
$ws.Range("H89").Value = 21278756
$ws.Range("I89").Value = 33335510
$ws.Range("J89").Value = 2127.7058
$ws.Range("K89").Value = 166677550
$ws.Range("L89").Value = 10638.529
$ws.Range("M89").Value = -166671934
$ws.Range("N89").Value = -21870.529
$ws.Range("H94").Value = 50393.168
$ws.Range("I94").Value = 432.9375
$ws.Range("J94").Value = 450075
$ws.Range("K94").Value = 432.9375
$ws.Range("L94").Value = 450075
$ws.Range("M94").Value = 18.0625
$ws.Range("N94").Value = -450977
$ws.Range("H99").Value = 1931.0294
$ws.Range("I99").Value = 2001.96
$ws.Range("K99").Value = 2001.96
$ws.Range("M99").Value = -503.96
$ws.Range("H134").Value = 2957.9285
$ws.Range("I134").Value = 2954.1924
$ws.Range("J134").Value = 3006.5
$ws.Range("K134").Value = 8862.5772
$ws.Range("L134").Value = 9019.5
$ws.Range("M134").Value = -6327.5772
$ws.Range("N134").Value = -14089.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2238.3572
$ws.Range("J31").Value = 2565.9783
$ws.Range("L31").Value = 2565.9783
$ws.Range("N31").Value = -3155.9783
$ws.Range("H34").Value = 2238.3572
$ws.Range("J34").Value = 2565.9783
$ws.Range("L34").Value = 2565.9783
$ws.Range("N34").Value = -2969.9783
$ws.Range("H58").Value = 2138.361
$ws.Range("I58").Value = 1928.2273
$ws.Range("J58").Value = 2468.5715
$ws.Range("K58").Value = 1928.2273
$ws.Range("L58").Value = 2468.5715
$ws.Range("M58").Value = -1725.2273
$ws.Range("N58").Value = -2874.5715
$ws.Range("H86").Value = 11977
$ws.Range("I86").Value = 9806.091
$ws.Range("J86").Value = 14962
$ws.Range("K86").Value = 9806.091
$ws.Range("L86").Value = 14962
$ws.Range("M86").Value = -8683.091
$ws.Range("N86").Value = -17208
$ws.Range("H89").Value = 11977
$ws.Range("I89").Value = 9806.091
$ws.Range("J89").Value = 14962
$ws.Range("K89").Value = 49030.455
$ws.Range("L89").Value = 74810
$ws.Range("M89").Value = -43414.455
$ws.Range("N89").Value = -86042
$ws.Range("H94").Value = 1301
$ws.Range("I94").Value = 1475.5834
$ws.Range("J94").Value = 1126.4166
$ws.Range("K94").Value = 1475.5834
$ws.Range("L94").Value = 1126.4166
$ws.Range("M94").Value = -1024.5834
$ws.Range("N94").Value = -2028.4166
$ws.Range("H96").Value = 16078
$ws.Range("J96").Value = 16078
$ws.Range("L96").Value = 16078
$ws.Range("N96").Value = -21570
$ws.Range("H132").Value = 2605.5557
$ws.Range("I132").Value = 2170.5881
$ws.Range("K132").Value = 6511.7643
$ws.Range("M132").Value = -3981.7643
$ws.Range("H136").Value = 2138.361
$ws.Range("I136").Value = 1928.2273
$ws.Range("J136").Value = 2468.5715
$ws.Range("K136").Value = 5784.6819
$ws.Range("L136").Value = 7405.7145
$ws.Range("M136").Value = -3234.6819
$ws.Range("N136").Value = -12505.7145

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 118158.37
$ws.Range("J37").Value = 118158.37
$ws.Range("L37").Value = 354475.11
$ws.Range("N37").Value = -354699.11
$ws.Range("H70").Value = 3858.8572
$ws.Range("I70").Value = 3012
$ws.Range("K70").Value = 9036
$ws.Range("M70").Value = -8721
$ws.Range("H73").Value = 3858.8572
$ws.Range("I73").Value = 3012
$ws.Range("K73").Value = 9036
$ws.Range("M73").Value = -7944
$ws.Range("H117").Value = 1727175
$ws.Range("I117").Value = 2900
$ws.Range("J117").Value = 6900000
$ws.Range("K117").Value = 8700
$ws.Range("L117").Value = 20700000
$ws.Range("M117").Value = -5258
$ws.Range("N117").Value = -20706884

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 107199.8
$ws.Range("J95").Value = 107199.8
$ws.Range("L95").Value = 107199.8
$ws.Range("N95").Value = -112691.8
$ws.Range("H102").Value = 49811.668
$ws.Range("I102").Value = 87024.586
$ws.Range("J102").Value = 12598.75
$ws.Range("K102").Value = 87024.586
$ws.Range("L102").Value = 12598.75
$ws.Range("M102").Value = -85402.586
$ws.Range("N102").Value = -15842.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1081.4736
$ws.Range("I82").Value = 1020.0909
$ws.Range("K82").Value = 1020.0909
$ws.Range("M82").Value = -659.0909
$ws.Range("H85").Value = 1081.4736
$ws.Range("I85").Value = 1020.0909
$ws.Range("K85").Value = 1020.0909
$ws.Range("M85").Value = 227.9091
$ws.Range("H93").Value = 2578.976
$ws.Range("I93").Value = 2263.5925
$ws.Range("J93").Value = 3146.6667
$ws.Range("K93").Value = 2263.5925
$ws.Range("L93").Value = 3146.6667
$ws.Range("M93").Value = -1015.5925
$ws.Range("N93").Value = -5642.6667
$ws.Range("H132").Value = 4214.5
$ws.Range("I132").Value = 3182.5454
$ws.Range("K132").Value = 9547.636200000001
$ws.Range("M132").Value = -7017.636200000001
$ws.Range("H136").Value = 41849.305
$ws.Range("I136").Value = 2663
$ws.Range("K136").Value = 7989
$ws.Range("M136").Value = -5439
